$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.645.92"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.398.32"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "'561.76"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "2.403.87"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "'26.04"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "2.828.83"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "60.497.14"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "2.408.95"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'8.07"
$ws.Range("E19").Value = "  +6.91%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'323.20"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'6.06"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'560.34"
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = "  -5.10%  "
$ws.Range("D29").Value = "2.516.40"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "0.0₃0934"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'18.26"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'41.75"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("E45").Value = "  +5.19%  "
$ws.Range("D46").Value = "0.0₆0275"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").Value = "'141.58"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "'3.52"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "'0.0504"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "'19.25"
$ws.Range("E51").Value = "  -1.69%  "
